$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H
$ws.Range("H1").Value = "CLR"

# Repeating r/g/b pattern for rows 2..22
$colors = @("r", "g", "b")
for ($row = 2; $row -le 22; $row++) {
    $idx = ($row - 2) % 3
    $ws.Cells.Item($row, 8).Value = $colors[$idx]
}

# Apply the same centered style as the rest of the table (style index 1:
# horizontal=center vertical=center) to the new column H cells.
$ws.Range("H1:H22").HorizontalAlignment = -4108
$ws.Range("H1:H22").VerticalAlignment = -4108

# Update the view: move the active selection to G11.
$ws.Range("G11").Select()
